$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern used below for cells whose new value looks like a plain
# decimal number (single '.') so that Excel's COM layer would otherwise
# auto-convert the text into a Number cell. We briefly mark the cell as
# Text ("@"), assign the literal string, then clear the number format
# again so the cell keeps using the default/general style (matching the
# original workbook, which stores these as plain text cells with no
# explicit style).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "34.013.72"
$ws.Range("E2").Value = "  -1.77%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.780.44"
$ws.Range("E3").Value = "  -3.99%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.37%  "

# Row 5 - BNB
Set-TextValue "D5" "223.76"
$ws.Range("E5").Value = "  -1.40%  "

# Row 6 - XRP
Set-TextValue "D6" "0.546"
$ws.Range("E6").Value = "  -1.97%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.38%  "

# Row 8 - Solana
$ws.Range("E8").Value = "  -1.20%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.282"
$ws.Range("E9").Value = "  -4.50%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.47%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.44%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "2.036.42"
$ws.Range("E12").Value = "  -3.77%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.782.20"
$ws.Range("E13").Value = "  -4.11%  "

# Row 14 - Chainlink
Set-TextValue "D14" "10.75"
$ws.Range("E14").Value = "  -4.50%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "33.984.30"
$ws.Range("E15").Value = "  -1.89%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  -5.50%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -5.79%  "

# Row 18 - Litecoin
Set-TextValue "D18" "67.44"
$ws.Range("E18").Value = "  -3.82%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "241.83"
$ws.Range("E19").Value = "  -4.99%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -3.91%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.30%  "

# Row 22 - Avalanche
Set-TextValue "D22" "10.64"
$ws.Range("E22").Value = "  -6.25%  "

# Row 23 - Uniswap
Set-TextValue "D23" "4.07"
$ws.Range("E23").Value = "  -6.49%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.83%  "

# Row 25 - Monero
Set-TextValue "D25" "159.47"
$ws.Range("E25").Value = "  -1.45%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "16.23"
$ws.Range("E26").Value = "  -4.54%  "

# Row 27 - Cosmos
Set-TextValue "D27" "7.00"
$ws.Range("E27").Value = "  -4.06%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -3.44%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.42%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -4.50%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.22%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -5.19%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "3.48"
$ws.Range("E33").Value = "  -4.88%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -8.40%  "

# Row 35 - Maker
Set-TextValue "D35" "1.390.30"
$ws.Range("E35").Value = "  -4.36%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "0.640"
$ws.Range("E36").Value = "  -3.58%  "

# Row 37 - TrustWalletToken
Set-TextValue "D37" "1.05"
$ws.Range("E37").Value = "  -3.03%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.0184"
$ws.Range("E38").Value = "  -5.63%  "

# Row 39 - was HuobiToken, becomes RenderToken
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D39" "2.21"
$ws.Range("E39").Value = "  +1.67%  "

# Row 40 - was RenderToken, becomes HuobiToken
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D40" "2.35"
$ws.Range("E40").Value = "  -0.74%  "

# Row 41 - MXToken
$ws.Range("E41").Value = "  -6.44%  "

# Row 42 - Aave
Set-TextValue "D42" "78.05"
$ws.Range("E42").Value = "  -6.73%  "

# Row 43 - ARBITRUM
Set-TextValue "D43" "0.907"
$ws.Range("E43").Value = "  -8.45%  "

# Row 44 - BabyDogeCoin (contains a literal subscript-six character, U+2086)
$ws.Range("D44").Value = "0.0₆0144"
$ws.Range("E44").Value = "  +12.99%  "

# Row 45 - WEMIXToken
$ws.Range("E45").Value = "  +1.46%  "

# Row 46 - Kaspa
$ws.Range("E46").Value = "  +0.53%  "

# Row 47 - Quant
Set-TextValue "D47" "106.68"
$ws.Range("E47").Value = "  -0.21%  "

# Row 48 - FraxShare
Set-TextValue "D48" "5.83"
$ws.Range("E48").Value = "  -4.88%  "

# Row 49 - RocketPoolETH
Set-TextValue "D49" "1.937.14"
$ws.Range("E49").Value = "  -3.36%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "12.15"
$ws.Range("E50").Value = "  -3.70%  "

# Row 51 - PaxDollar
$ws.Range("E51").Value = "  -0.18%  "
